# "palabras filtro: rede de monitoeo y calidad del agua"
# Replace the ECCA group row with the "Gestión Integrada de Recursos Hídricos"
# group and its single working-paper entry, and drop the other now-unused
# group rows (Oceanografía Operacional, IDENTUS, GICMIL).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 4-6 entirely (delete bottom-up so row indices stay valid).
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Update row 3 with the new group name and its document entry.
$ws.Range("A3").Value = "Gestión Integrada de Recursos Hídricos"
$ws.Range("B3").Value = "8.- Documento de trabajo (Working Paper) : Red de monitoreo de la calidad del agua del río cauca y sus tributarios. Tramos Salvajina - La Virginia. Vol. VI. Fase III  2007,  Nro. Paginas: 103,  Instituciones participantes: Corporación Autónoma Regional del Valle del Cauca, CVC Universidad del Valle Facultad de Ingeniería Escuela de Ingeniería de Recursos Naturales y del Ambiente, EIDENAR Instituto Cinara,  URL: http://www.cvc.gov.co/cvc/Mosaic/,  DOI: Autores: CARLOS ALBERTO RAMIREZ CALLEJAS, ALBERTO GALVIS CASTANO, DIANA PAOLA BERNAL SUAREZ"
